$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Unprotect()

# Update the confidential disclaimer cell (A44) with the new "as of" date.
$ws.Range("A44").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-23 for illustrative purposes only and are subject to change."

# Update the Weight (D) and Percent Change (E) values for each holding row.
$updates = @(
    @{ Row = 2; D = 0.01908009950491013; E = -0.008067835679591795 },
    @{ Row = 3; D = 0.0166344162144561; E = -0.01000944287063255 },
    @{ Row = 4; D = 0.007106304286668637; E = 0.0172188943657301 },
    @{ Row = 5; D = 0.03189739628190167; E = 0.005239598363117937 },
    @{ Row = 6; D = 0.0488645008206469; E = 0.00856030628087967 },
    @{ Row = 7; D = 0.07752674661762941; E = -0.006888726801199452 },
    @{ Row = 8; D = 0.01978505946685988; E = 0.0003334444814939719 },
    @{ Row = 9; D = 0.03016908200798608; E = -0.02018056293149229 },
    @{ Row = 10; D = 0.05112341998570916; E = -0.008188208979070133 },
    @{ Row = 11; D = 0.007465538563179128; E = -0.02360713684564886 },
    @{ Row = 12; D = 0.01771730257311957; E = -0.005558806319485021 },
    @{ Row = 13; D = 0.01896229202463778; E = -0.005765407554671809 },
    @{ Row = 14; D = 0.01922618078044785; E = 0.007647058823529562 },
    @{ Row = 15; D = 0.02213602554317494; E = -0.0040801844952989 },
    @{ Row = 16; D = 0.02107717191048704; E = -0.01234126274369518 },
    @{ Row = 17; D = 0.03227579390853647; E = -0.009913470055188611 },
    @{ Row = 18; D = 0.02723143467991471; E = 0.003599363189589644 },
    @{ Row = 19; D = 0.02371056178617502; E = -0.03275940880694794 },
    @{ Row = 20; D = 0.037816201167425; E = -0.0009345794392523477 },
    @{ Row = 21; D = 0.04268494871212076; E = -0.01000198714976475 },
    @{ Row = 22; D = 0.02814813395307398; E = -0.008733258928571463 },
    @{ Row = 23; D = 0.01826251559182001; E = -0.0159979357502259 },
    @{ Row = 24; D = 0.06672332944673363; E = 0.006737573626001359 },
    @{ Row = 25; D = 0.008216836134036008; E = 0.02290149299382538 },
    @{ Row = 26; D = 0.01954708835670974; E = -0.008317060823027522 },
    @{ Row = 27; D = 0.008285007395953611; E = -0.008759124087591275 },
    @{ Row = 28; D = 0.01666771646221309; E = 0.01432448733413749 },
    @{ Row = 29; D = 0.007677434951028999; E = -0.002516521063076649 },
    @{ Row = 30; D = 0.01731612883363211; E = 0.01204644412191569 },
    @{ Row = 31; D = 0.01809522896983327; E = -0.01777777777777778 },
    @{ Row = 32; D = 0.02862753185946227; E = 0.01858964510677508 },
    @{ Row = 33; D = 0.006759636141387195; E = -0.0007900729655621053 },
    @{ Row = 34; D = 0.02104826980866023; E = -0.01170149253731345 },
    @{ Row = 35; D = 0.007141646530750343; E = -0.005828531210135157 },
    @{ Row = 36; D = 0.02881256814181004; E = 0.004470370168456528 },
    @{ Row = 37; D = 0.02133352098755968; E = 0.005654709312599371 },
    @{ Row = 38; D = 0.03920632943463875; E = 0.0007211538461540101 },
    @{ Row = 39; D = 0.03534679930427614; E = 0.01186069351061425 },
    @{ Row = 40; D = 0.03029380086043441; E = -0.02141449756299918 },
    @{ Row = 41; D = 0.9999999999999997; E = -0.002957910214678283 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 4).Value = $u.D
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
